# Auto-generated Excel COM-interop script to update the cryptos worksheet
# (GitHub Actions daily refresh of coinranking.com price/volume snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store free-form text that can LOOK
# numeric (e.g. "27.319.15", "0.9951", "  +9.46%  "). Excel's default
# Value setter auto-coerces such strings into real numbers, which both
# changes the cell type and mangles the exact decimal text (float drift,
# percent signs stripped, etc). Force Text format first on exactly the
# cells being rewritten so the literal source strings are preserved.
# (Kept as separate contiguous-range calls rather than one comma-joined
# multi-area range, since multi-area NumberFormat assignment only takes
# effect on the first area in this host.)
$ws.Range("D2:D25").NumberFormat = "@"
$ws.Range("D27:D42").NumberFormat = "@"
$ws.Range("D44:D51").NumberFormat = "@"
$ws.Range("E2:E48").NumberFormat = "@"
$ws.Range("E50:E51").NumberFormat = "@"

# Rows with Price (D) and Volume(1h) (E) updates
$ws.Range("D2").Value = '27.434.65'
$ws.Range("E2").Value = '  +9.89%  '
$ws.Range("D3").Value = '1.763.66'
$ws.Range("E3").Value = '  +5.39%  '
$ws.Range("D4").Value = '0.9951'
$ws.Range("E4").Value = '  -0.68%  '
$ws.Range("D5").Value = '335.37'
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").Value = '0.9912'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '0.3768'
$ws.Range("E7").Value = '  +3.03%  '
$ws.Range("D8").Value = '48.88'
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("D9").Value = '0.3436'
$ws.Range("E9").Value = '  +6.40%  '
$ws.Range("D10").Value = '1.207'
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("D11").Value = '0.07615'
$ws.Range("E11").Value = '  +6.46%  '
$ws.Range("D12").Value = '0.9903'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '6.482'
$ws.Range("E13").Value = '  +6.43%  '
$ws.Range("D14").Value = '21.08'
$ws.Range("E14").Value = '  +7.24%  '
$ws.Range("D15").Value = '7.135'
$ws.Range("E15").Value = '  +7.17%  '
$ws.Range("D16").Value = '1.758.27'
$ws.Range("E16").Value = '  +5.56%  '
$ws.Range("D17").Value = '0.00001102'
$ws.Range("E17").Value = '  +5.08%  '
$ws.Range("D18").Value = '0.06736'
$ws.Range("E18").Value = '  +2.80%  '
$ws.Range("D19").Value = '84.16'
$ws.Range("E19").Value = '  +6.83%  '
$ws.Range("D20").Value = '0.9921'
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").Value = '17.22'
$ws.Range("E21").Value = '  +8.67%  '
$ws.Range("D22").Value = '6.296'
$ws.Range("E22").Value = '  +6.41%  '
$ws.Range("D23").Value = '27.392.67'
$ws.Range("E23").Value = '  +9.77%  '
$ws.Range("D24").Value = '13.01'
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").Value = '2.455'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D27").Value = '2.467'
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("D28").Value = '153.06'
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("D29").Value = '19.93'
$ws.Range("E29").Value = '  +6.65%  '
$ws.Range("D30").Value = '1.955.68'
$ws.Range("E30").Value = '  +5.42%  '
$ws.Range("D31").Value = '134.30'
$ws.Range("E31").Value = '  +6.66%  '
$ws.Range("D32").Value = '4.093'
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").Value = '6.161'
$ws.Range("E33").Value = '  +6.28%  '
$ws.Range("D34").Value = '0.08647'
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("D37").Value = '5.517'
$ws.Range("E37").Value = '  +6.83%  '
$ws.Range("D38").Value = '0.02381'
$ws.Range("E38").Value = '  +6.72%  '
$ws.Range("D39").Value = '0.06391'
$ws.Range("E39").Value = '  +5.53%  '
$ws.Range("D40").Value = '0.2214'
$ws.Range("E40").Value = '  +5.80%  '
$ws.Range("D44").Value = '14.37'
$ws.Range("E44").Value = '  +4.80%  '
$ws.Range("D45").Value = '0.9919'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").Value = '0.6342'
$ws.Range("E46").Value = '  +10.67%  '
$ws.Range("D47").Value = '3.938'
$ws.Range("E47").Value = '  +2.39%  '
$ws.Range("D48").Value = '2.113'
$ws.Range("E48").Value = '  +7.67%  '
$ws.Range("D50").Value = '0.07305'
$ws.Range("E50").Value = '  +4.31%  '
$ws.Range("D51").Value = '79.18'
$ws.Range("E51").Value = '  +6.04%  '

# Rows with only Volume(1h) (E) updates
$ws.Range("E26").Value = '  +26.18%  '
$ws.Range("E43").Value = '  -0.02%  '

# Rows with only Price (D) updates
$ws.Range("D49").Value = '130.37'

# Rows 35/36 and 41/42: the coin ranking order shifted, so the Coin name,
# Link, Price and Volume all changed for these row positions.
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '13.10'
$ws.Range("E35").Value = '  +6.32%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.693'
$ws.Range("E36").Value = '  +1.56%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6579'
$ws.Range("E41").Value = '  +10.34%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.633'
$ws.Range("E42").Value = '  +4.69%  '
